$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting for the brand-new rows (30:43) ---
# Column A: reuse the bold/bordered label style already used throughout
# column A (copy format only from an existing labeled cell).
$ws.Range("A2").Copy()
$ws.Range("A30:A43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column B: reuse the existing "n_sample" text cells (B2:B8 = 0,100,500,
# 1000,2500,5000,"5426 (all)") so the new rows get the same text-typed
# values (not auto-converted to numbers) with matching formatting.
$ws.Range("B2:B8").Copy()
$ws.Range("B30").PasteSpecial(-4104)
$ws.Range("B37").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Row 16: accuracy_balanced_mean
$ws.Cells.Item(16,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = 0.393

# Row 17: accuracy_balanced_mean
$ws.Cells.Item(17,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(17,3).Value = 0.148
$ws.Cells.Item(17,4).Value = 0.121
$ws.Cells.Item(17,5).Value = 0.238
$ws.Cells.Item(17,6).Value = 0.204
$ws.Cells.Item(17,7).Value = 0.189
$ws.Cells.Item(17,8).Value = 0.353

# Row 18: accuracy_balanced_mean
$ws.Cells.Item(18,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(18,3).Value = 0.333
$ws.Cells.Item(18,4).Value = 0.305
$ws.Cells.Item(18,5).Value = 0.4
$ws.Cells.Item(18,6).Value = 0.358
$ws.Cells.Item(18,7).Value = 0.507
$ws.Cells.Item(18,8).Value = 0.576

# Row 19: accuracy_balanced_mean
$ws.Cells.Item(19,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(19,3).Value = 0.392
$ws.Cells.Item(19,4).Value = 0.374
$ws.Cells.Item(19,5).Value = 0.451
$ws.Cells.Item(19,6).Value = 0.409
$ws.Cells.Item(19,7).Value = 0.534
$ws.Cells.Item(19,8).Value = 0.606

# Row 20: accuracy_balanced_mean
$ws.Cells.Item(20,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(20,3).Value = 0.448
$ws.Cells.Item(20,4).Value = 0.463
$ws.Cells.Item(20,5).Value = 0.472
$ws.Cells.Item(20,6).Value = 0.488
$ws.Cells.Item(20,7).Value = 0.574
$ws.Cells.Item(20,8).Value = 0.638

# Row 21: accuracy_balanced_mean
$ws.Cells.Item(21,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(21,3).Value = 0.489
$ws.Cells.Item(21,4).Value = 0.515
$ws.Cells.Item(21,5).Value = 0.499
$ws.Cells.Item(21,6).Value = 0.505
$ws.Cells.Item(21,7).Value = 0.609
$ws.Cells.Item(21,8).Value = 0.649

# Row 22: accuracy_balanced_mean
$ws.Cells.Item(22,1).Value = "accuracy_balanced_mean"
$ws.Cells.Item(22,3).Value = 0.478
$ws.Cells.Item(22,4).Value = 0.516
$ws.Cells.Item(22,5).Value = 0.507
$ws.Cells.Item(22,6).Value = 0.516
$ws.Cells.Item(22,7).Value = 0.596
$ws.Cells.Item(22,8).Value = 0.657

# Row 23: f1_macro_std
$ws.Cells.Item(23,1).Value = "f1_macro_std"
$ws.Cells.Item(23,3).Value = 0
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 0
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = 0

# Row 24: f1_macro_std
$ws.Cells.Item(24,1).Value = "f1_macro_std"
$ws.Cells.Item(24,3).Value = 0.018
$ws.Cells.Item(24,4).Value = 0.012
$ws.Cells.Item(24,5).Value = 0.021
$ws.Cells.Item(24,6).Value = 0.031
$ws.Cells.Item(24,7).Value = 0.105
$ws.Cells.Item(24,8).Value = 0.126

# Row 25: f1_macro_std
$ws.Cells.Item(25,1).Value = "f1_macro_std"
$ws.Cells.Item(25,3).Value = 0.013
$ws.Cells.Item(25,4).Value = 0.027
$ws.Cells.Item(25,5).Value = 0.013
$ws.Cells.Item(25,6).Value = 0.011
$ws.Cells.Item(25,7).Value = 0.016
$ws.Cells.Item(25,8).Value = 0.015

# Row 26: f1_macro_std
$ws.Cells.Item(26,1).Value = "f1_macro_std"
$ws.Cells.Item(26,3).Value = 0.009
$ws.Cells.Item(26,4).Value = 0.008
$ws.Cells.Item(26,5).Value = 0.005
$ws.Cells.Item(26,6).Value = 0.006
$ws.Cells.Item(26,7).Value = 0.009
$ws.Cells.Item(26,8).Value = 0.011

# Row 27: f1_macro_std
$ws.Cells.Item(27,1).Value = "f1_macro_std"
$ws.Cells.Item(27,3).Value = 0.012
$ws.Cells.Item(27,4).Value = 0.006
$ws.Cells.Item(27,5).Value = 0.015
$ws.Cells.Item(27,6).Value = 0.019
$ws.Cells.Item(27,7).Value = 0.01
$ws.Cells.Item(27,8).Value = 0.01

# Row 28: f1_macro_std
$ws.Cells.Item(28,1).Value = "f1_macro_std"
$ws.Cells.Item(28,3).Value = 0.006
$ws.Cells.Item(28,4).Value = 0.006
$ws.Cells.Item(28,5).Value = 0.006
$ws.Cells.Item(28,6).Value = 0.006
$ws.Cells.Item(28,7).Value = 0.008
$ws.Cells.Item(28,8).Value = 0.004

# Row 29: f1_macro_std
$ws.Cells.Item(29,1).Value = "f1_macro_std"
$ws.Cells.Item(29,3).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = 0.002
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 0.009
$ws.Cells.Item(29,8).Value = 0.011

# Row 30: f1_micro_std
$ws.Cells.Item(30,1).Value = "f1_micro_std"
$ws.Cells.Item(30,3).Value = 0
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 0
$ws.Cells.Item(30,8).Value = 0

# Row 31: f1_micro_std
$ws.Cells.Item(31,1).Value = "f1_micro_std"
$ws.Cells.Item(31,3).Value = 0.018
$ws.Cells.Item(31,4).Value = 0.007
$ws.Cells.Item(31,5).Value = 0.016
$ws.Cells.Item(31,6).Value = 0.026
$ws.Cells.Item(31,7).Value = 0.094
$ws.Cells.Item(31,8).Value = 0.159

# Row 32: f1_micro_std
$ws.Cells.Item(32,1).Value = "f1_micro_std"
$ws.Cells.Item(32,3).Value = 0.004
$ws.Cells.Item(32,4).Value = 0.007
$ws.Cells.Item(32,5).Value = 0.016
$ws.Cells.Item(32,6).Value = 0.011
$ws.Cells.Item(32,7).Value = 0.008
$ws.Cells.Item(32,8).Value = 0.006

# Row 33: f1_micro_std
$ws.Cells.Item(33,1).Value = "f1_micro_std"
$ws.Cells.Item(33,3).Value = 0.001
$ws.Cells.Item(33,4).Value = 0.005
$ws.Cells.Item(33,5).Value = 0.01
$ws.Cells.Item(33,6).Value = 0.008
$ws.Cells.Item(33,7).Value = 0.011
$ws.Cells.Item(33,8).Value = 0.012

# Row 34: f1_micro_std
$ws.Cells.Item(34,1).Value = "f1_micro_std"
$ws.Cells.Item(34,3).Value = 0.004
$ws.Cells.Item(34,4).Value = 0.003
$ws.Cells.Item(34,5).Value = 0.011
$ws.Cells.Item(34,6).Value = 0.005
$ws.Cells.Item(34,7).Value = 0.004
$ws.Cells.Item(34,8).Value = 0.005

# Row 35: f1_micro_std
$ws.Cells.Item(35,1).Value = "f1_micro_std"
$ws.Cells.Item(35,3).Value = 0.002
$ws.Cells.Item(35,4).Value = 0.005
$ws.Cells.Item(35,5).Value = 0.001
$ws.Cells.Item(35,6).Value = 0.003
$ws.Cells.Item(35,7).Value = 0.007
$ws.Cells.Item(35,8).Value = 0.012

# Row 36: f1_micro_std
$ws.Cells.Item(36,1).Value = "f1_micro_std"
$ws.Cells.Item(36,3).Value = 0
$ws.Cells.Item(36,4).Value = 0
$ws.Cells.Item(36,5).Value = 0
$ws.Cells.Item(36,6).Value = 0
$ws.Cells.Item(36,7).Value = 0.004
$ws.Cells.Item(36,8).Value = 0.007

# Row 37: accuracy_balanced_std
$ws.Cells.Item(37,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(37,3).Value = 0
$ws.Cells.Item(37,4).Value = 0
$ws.Cells.Item(37,5).Value = 0
$ws.Cells.Item(37,6).Value = 0
$ws.Cells.Item(37,7).Value = 0
$ws.Cells.Item(37,8).Value = 0

# Row 38: accuracy_balanced_std
$ws.Cells.Item(38,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(38,3).Value = 0.015
$ws.Cells.Item(38,4).Value = 0.007
$ws.Cells.Item(38,5).Value = 0.024
$ws.Cells.Item(38,6).Value = 0.032
$ws.Cells.Item(38,7).Value = 0.112
$ws.Cells.Item(38,8).Value = 0.159

# Row 39: accuracy_balanced_std
$ws.Cells.Item(39,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(39,3).Value = 0.013
$ws.Cells.Item(39,4).Value = 0.019
$ws.Cells.Item(39,5).Value = 0.009
$ws.Cells.Item(39,6).Value = 0.011
$ws.Cells.Item(39,7).Value = 0.022
$ws.Cells.Item(39,8).Value = 0.011

# Row 40: accuracy_balanced_std
$ws.Cells.Item(40,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(40,3).Value = 0.008
$ws.Cells.Item(40,4).Value = 0.011
$ws.Cells.Item(40,5).Value = 0.008
$ws.Cells.Item(40,6).Value = 0.004
$ws.Cells.Item(40,7).Value = 0.005
$ws.Cells.Item(40,8).Value = 0.034

# Row 41: accuracy_balanced_std
$ws.Cells.Item(41,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(41,3).Value = 0.014
$ws.Cells.Item(41,4).Value = 0.008
$ws.Cells.Item(41,5).Value = 0.014
$ws.Cells.Item(41,6).Value = 0.021
$ws.Cells.Item(41,7).Value = 0.019
$ws.Cells.Item(41,8).Value = 0.012

# Row 42: accuracy_balanced_std
$ws.Cells.Item(42,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(42,3).Value = 0.004
$ws.Cells.Item(42,4).Value = 0.004
$ws.Cells.Item(42,5).Value = 0.006
$ws.Cells.Item(42,6).Value = 0.005
$ws.Cells.Item(42,7).Value = 0.008
$ws.Cells.Item(42,8).Value = 0.008

# Row 43: accuracy_balanced_std
$ws.Cells.Item(43,1).Value = "accuracy_balanced_std"
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,5).Value = 0
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 0.013
$ws.Cells.Item(43,8).Value = 0.003
